# "Add validation for formatting"
#
# Adds 18 rows of repeated sample data (Israel / Jerusalem / 9.2 /
# 24-02-2021 / 22145) to the "Country Population" sheet, formatted with a
# Verdana 11pt font (mirrors the new font/fill/style entries that appear in
# the target workbook), and tweaks two header/format-flavoured cells
# (D1 -> right aligned, B2 -> filled) the same way the source edit does.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Country Population")

# --- formatting touch-ups on existing header/data cells -------------------
$ws.Range("D1").HorizontalAlignment = -4152          # xlRight
$ws.Range("B2").Interior.ColorIndex = 11              # flag a fill on B2

# --- 18 new sample rows (6-23), each holding the same 5 values -----------
$country = "Israel"
$capital = "Jerusalem"
$population = "9.2"
$surveyDate = "24-02-2021"
$extra = "22145"

for ($i = 0; $i -lt 18; $i++) {
    $row = 6 + $i

    $ws.Cells.Item($row, 1).Value = $country
    $ws.Cells.Item($row, 2).Value = $capital

    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = $population

    $ws.Cells.Item($row, 4).NumberFormat = "@"
    $ws.Cells.Item($row, 4).Value = $surveyDate

    $ws.Cells.Item($row, 5).NumberFormat = "@"
    $ws.Cells.Item($row, 5).Value = $extra

    $rowRange = $ws.Range("A" + $row + ":E" + $row)
    $rowRange.Font.Name = "Verdana"
    $rowRange.Font.Size = 11
}
